$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data, continuing the trade log pattern used by rows 2-4.
$ws.Cells.Item(5, 1).Value = 42647.680543981478
$ws.Cells.Item(5, 2).Value = $true
$ws.Cells.Item(5, 3).Value = 10102.14
$ws.Cells.Item(5, 4).Value = 10039.89
$ws.Cells.Item(5, 5).Value = 78.63
$ws.Cells.Item(5, 6).Value = 78.14
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = -0.62
$ws.Cells.Item(5, 9).Value = $true

# Match the date-time formatting used for the Date and IsShortSell columns
# (A3/G3) by copying their formats onto the new row's cells, so the same
# shared style entry is reused instead of creating a new one.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)

$ws.Cells.Item(3, 7).Copy()
$ws.Cells.Item(5, 7).PasteSpecial(-4122)

$excel.CutCopyMode = $false
